$wb = $excel.ActiveWorkbook

# --- Hoja1!A1 : update the daily conversion message with new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$a1 = $ws1.Range("A1")
$text = $a1.Value()
$text = $text.Replace("1000 Bs = 3.14 = 11871.26 pesos", "1000 Bs = 3.14 = 11825.51 pesos")
$text = $text.Replace("11871.26 pesos = 3.13 = 971.63 Bs", "11825.51 pesos = 3.12 = 968.28 Bs")
$a1.Value = $text

# --- tasas sheet : updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 3760.5
$ws2.Range("N12").Value = 3786
$ws2.Range("O12").Value = 310
